$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = (Get-Date -Year 2018 -Month 3 -Day 4).Date
$ws.Range("B9").Value = "Amélioration du cahier des charges et du Gantt"
$ws.Range("C9").Value = 1

$ws.Range("A10").Value = (Get-Date -Year 2018 -Month 3 -Day 5).Date
$ws.Range("B10").Value = "Amélioration du cahier des charges et du Gantt"
$ws.Range("C10").Value = 4

$ws.Range("C11").Select()
